# Generate Report for Handoff
# 189213fe-347a-4762-9be1-2ae0b57be902 has now also been handed off (alongside
# 8f54520d-3dac-487a-81ea-6aa7df61bbbe which already was), so both files move
# to "Ready for handoff" status with a freshly generated handoff timestamp,
# and the report (sorted by file name) now lists 189213fe... before 8f54520d....

$wb = $excel.ActiveWorkbook

function Set-RowNineTen-Overview($sheetName) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Row 9 becomes 189213fe..., Row 10 becomes 8f54520d..., both "Ready for handoff"
    $ws.Range("A9").Value  = "189213fe-347a-4762-9be1-2ae0b57be902.md"
    $ws.Range("B9").Value  = "Ready for handoff"
    $ws.Range("C9").Value  = "Ready for handoff"

    $ws.Range("A10").Value = "8f54520d-3dac-487a-81ea-6aa7df61bbbe.md"
    $ws.Range("B10").Value = "Ready for handoff"
    $ws.Range("C10").Value = "Ready for handoff"

    foreach ($h in $ws.Hyperlinks) {
        $addr = $h.Range.Address()
        if ($addr -eq '$A$9') {
            $h.TextToDisplay = "189213fe-347a-4762-9be1-2ae0b57be902.md"
        }
        elseif ($addr -eq '$A$10') {
            $h.TextToDisplay = "8f54520d-3dac-487a-81ea-6aa7df61bbbe.md"
        }
    }
}

function Set-RowNineTen-Locale($sheetName, $locale, $handoffDatetime) {
    $ws = $wb.Worksheets.Item($sheetName)

    $newXlf9  = "189213fe-347a-4762-9be1-2ae0b57be902.f4c13c59ebad7e5e1c38915a39fd2517f81dfdd8.$locale.xlf"
    $newXlf10 = "8f54520d-3dac-487a-81ea-6aa7df61bbbe.9335be45d1b69b758e4eeb5c62a91445b3952c35.$locale.xlf"

    # Row 9 becomes 189213fe..., Row 10 becomes 8f54520d..., both "Ready for handoff"
    $ws.Range("A9").Value = "189213fe-347a-4762-9be1-2ae0b57be902.md"
    $ws.Range("B9").Value = "Ready for handoff"
    $ws.Range("C9").Value = $newXlf9
    $ws.Range("D9").Value = $handoffDatetime

    $ws.Range("A10").Value = "8f54520d-3dac-487a-81ea-6aa7df61bbbe.md"
    $ws.Range("B10").Value = "Ready for handoff"
    $ws.Range("C10").Value = $newXlf10
    $ws.Range("D10").Value = $handoffDatetime

    foreach ($h in $ws.Hyperlinks) {
        $addr = $h.Range.Address()
        if ($addr -eq '$A$9') {
            $h.TextToDisplay = "189213fe-347a-4762-9be1-2ae0b57be902.md"
        }
        elseif ($addr -eq '$C$9') {
            $h.TextToDisplay = $newXlf9
        }
        elseif ($addr -eq '$A$10') {
            $h.TextToDisplay = "8f54520d-3dac-487a-81ea-6aa7df61bbbe.md"
        }
        elseif ($addr -eq '$C$10') {
            $h.TextToDisplay = $newXlf10
        }
    }
}

Set-RowNineTen-Overview "Overview"
Set-RowNineTen-Locale "zh-cn" "zh-cn" "2016-03-11 00:29:06"
Set-RowNineTen-Locale "de-de" "de-de" "2016-03-11 00:29:14"
